$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value2 = 0.04815801081838345
$ws.Range("D2").Value2 = 0.08424054524450497
$ws.Range("E2").Value2 = 0.04443085484498965
$ws.Range("F2").Value2 = 2.114851626990827
$ws.Range("G2").Value2 = 0.002494564884047172
$ws.Range("K2").Value2 = 2.204494428052158
$ws.Range("M2").Value2 = 0.5382924935363036
$ws.Range("N2").Value2 = 1.811455330370194
$ws.Range("C3").Value2 = 0.04278531419993215
$ws.Range("D3").Value2 = 0.08441390018681716
$ws.Range("E3").Value2 = 0.0410073458305682
$ws.Range("F3").Value2 = 2.053386877792008
$ws.Range("G3").Value2 = 0.002500996984734648
$ws.Range("K3").Value2 = 2.015053012493297
$ws.Range("M3").Value2 = 0.4929787337973153
$ws.Range("N3").Value2 = 1.822283627367696
$ws.Range("C4").Value2 = 0.03950802919261776
$ws.Range("D4").Value2 = 0.0845643176267501
$ws.Range("E4").Value2 = 0.03893700104591247
$ws.Range("F4").Value2 = 2.017091007264014
$ws.Range("G4").Value2 = 0.002505147741458703
$ws.Range("K4").Value2 = 1.900055044942292
$ws.Range("M4").Value2 = 0.4654915336702103
$ws.Range("N4").Value2 = 1.829572956835449
$ws.Range("C5").Value2 = 0.03817774873766666
$ws.Range("D5").Value2 = 0.08463657251753887
$ws.Range("E5").Value2 = 0.03810110082769214
$ws.Range("F5").Value2 = 2.002659378696222
$ws.Range("G5").Value2 = 0.002506890053935495
$ws.Range("K5").Value2 = 1.853518474428881
$ws.Range("M5").Value2 = 0.4543730635860825
$ws.Range("N5").Value2 = 1.832703532822933
$ws.Range("C6").Value2 = 0.03795716780361147
$ws.Range("D6").Value2 = 0.08464922936565245
$ws.Range("E6").Value2 = 0.03796276474302118
$ws.Range("F6").Value2 = 2.00028460430056
$ws.Range("G6").Value2 = 0.002507182440279953
$ws.Range("K6").Value2 = 1.845810635263092
$ws.Range("M6").Value2 = 0.4525318036908246
$ws.Range("N6").Value2 = 1.833233004149299
$ws.Range("C7").Value2 = 0.0394900675915153
$ws.Range("D7").Value2 = 0.08456524784125463
$ws.Range("E7").Value2 = 0.03892569654333755
$ws.Range("F7").Value2 = 2.016894927735677
$ws.Range("G7").Value2 = 0.002505171032754488
$ws.Range("K7").Value2 = 1.899426124023364
$ws.Range("M7").Value2 = 0.4653412530010357
$ws.Range("N7").Value2 = 1.829614529882669
$ws.Range("C8").Value2 = 0.04630091416186133
$ws.Range("D8").Value2 = 0.08429113008216405
$ws.Range("E8").Value2 = 0.04324373277054505
$ws.Range("F8").Value2 = 2.093356606572712
$ws.Range("G8").Value2 = 0.002496740986160235
$ws.Range("K8").Value2 = 2.138897489665908
$ws.Range("M8").Value2 = 0.5225976582091505
$ws.Range("N8").Value2 = 1.815055344888989
$ws.Range("C9").Value2 = 0.05983720648880819
$ws.Range("D9").Value2 = 0.08410712409052223
$ws.Range("E9").Value2 = 0.05197202544043478
$ws.Range("F9").Value2 = 2.254936329797005
$ws.Range("G9").Value2 = 0.002481798881615559
$ws.Range("K9").Value2 = 2.619264845396856
$ws.Range("M9").Value2 = 0.6376198606821362
$ws.Range("N9").Value2 = 1.791633774062859
$ws.Range("C10").Value2 = 0.06990584160188007
$ws.Range("D10").Value2 = 0.08419398476627293
$ws.Range("E10").Value2 = 0.0585569139738169
$ws.Range("F10").Value2 = 2.381014701380536
$ws.Range("G10").Value2 = 0.00247177709500929
$ws.Range("K10").Value2 = 2.97919999010287
$ws.Range("M10").Value2 = 0.7239186152899748
$ws.Range("N10").Value2 = 1.777615608614099
$ws.Range("C11").Value2 = 0.07451619938737508
$ws.Range("D11").Value2 = 0.08428316271416492
$ws.Range("E11").Value2 = 0.06159290872731304
$ws.Range("F11").Value2 = 2.440030815242295
$ws.Range("G11").Value2 = 0.002467422836409678
$ws.Range("K11").Value2 = 3.144566635144372
$ws.Range("M11").Value2 = 0.7635942300331067
$ws.Range("N11").Value2 = 1.771944835545114
$ws.Range("C12").Value2 = 0.07626658067854919
$ws.Range("D12").Value2 = 0.08432419396745416
$ws.Range("E12").Value2 = 0.06274861697752243
$ws.Range("F12").Value2 = 2.462622466457077
$ws.Range("G12").Value2 = 0.002465803218428819
$ws.Range("K12").Value2 = 3.207428742711159
$ws.Range("M12").Value2 = 0.7786805034048001
$ws.Range("N12").Value2 = 1.769900177398085
$ws.Range("C13").Value2 = 0.07588939945358675
$ws.Range("D13").Value2 = 0.08431503226285031
$ws.Range("E13").Value2 = 0.06249944222687986
$ws.Range("F13").Value2 = 2.457746046940599
$ws.Range("G13").Value2 = 0.002466150734501448
$ws.Range("K13").Value2 = 3.193879417403309
$ws.Range("M13").Value2 = 0.7754286174600509
$ws.Range("N13").Value2 = 1.770335942440525
$ws.Range("C14").Value2 = 0.07466011211205625
$ws.Range("D14").Value2 = 0.0842863921532313
$ws.Range("E14").Value2 = 0.061687867158561
$ws.Range("F14").Value2 = 2.441884537357538
$ws.Range("G14").Value2 = 0.002467289004222055
$ws.Range("K14").Value2 = 3.149733461460471
$ws.Range("M14").Value2 = 0.7648341337037863
$ws.Range("N14").Value2 = 1.771774554300251
$ws.Range("C15").Value2 = 0.07390773559311015
$ws.Range("D15").Value2 = 0.08426979851743965
$ws.Range("E15").Value2 = 0.06119154755114664
$ws.Range("F15").Value2 = 2.432200753735998
$ws.Range("G15").Value2 = 0.002467990030571039
$ws.Range("K15").Value2 = 3.122724427107983
$ws.Range("M15").Value2 = 0.7583528347610127
$ws.Range("N15").Value2 = 1.772669161265512
$ws.Range("C16").Value2 = 0.06960516735455258
$ws.Range("D16").Value2 = 0.08418916634264662
$ws.Range("E16").Value2 = 0.05835933814490346
$ws.Range("F16").Value2 = 2.377191692445734
$ws.Range("G16").Value2 = 0.002472065758684244
$ws.Range("K16").Value2 = 2.968426284396628
$ws.Range("M16").Value2 = 0.7213342898726012
$ws.Range("N16").Value2 = 1.778000518472339
$ws.Range("C17").Value2 = 0.06697353666454831
$ws.Range("D17").Value2 = 0.08415250589302303
$ws.Range("E17").Value2 = 0.05663239643077134
$ws.Range("F17").Value2 = 2.343874360340919
$ws.Range("G17").Value2 = 0.002474618379622405
$ws.Range("K17").Value2 = 2.874191805710836
$ws.Range("M17").Value2 = 0.698732971431852
$ws.Range("N17").Value2 = 1.781452842516501
$ws.Range("C18").Value2 = 0.06546271109304769
$ws.Range("D18").Value2 = 0.08413608713484422
$ws.Range("E18").Value2 = 0.05564290352611323
$ws.Range("F18").Value2 = 2.32486731533362
$ws.Range("G18").Value2 = 0.002476105858605401
$ws.Range("K18").Value2 = 2.820143658890402
$ws.Range("M18").Value2 = 0.6857725081482897
$ws.Range("N18").Value2 = 1.783504922538199
$ws.Range("C19").Value2 = 0.06495164923389041
$ws.Range("D19").Value2 = 0.08413132603512707
$ws.Range("E19").Value2 = 0.0553085239637312
$ws.Range("F19").Value2 = 2.318458549549888
$ws.Range("G19").Value2 = 0.002476612810441954
$ws.Range("K19").Value2 = 2.801869960745648
$ws.Range("M19").Value2 = 0.6813909966950291
$ws.Range("N19").Value2 = 1.784211088771116
$ws.Range("C20").Value2 = 0.06725338498232247
$ws.Range("D20").Value2 = 0.08415592454743859
$ws.Range("E20").Value2 = 0.05681583763692544
$ws.Range("F20").Value2 = 2.347404848643464
$ws.Range("G20").Value2 = 0.002474344654690108
$ws.Range("K20").Value2 = 2.884207338985789
$ws.Range("M20").Value2 = 0.701134847968774
$ws.Range("N20").Value2 = 1.781078456551683
$ws.Range("C21").Value2 = 0.07502105870403852
$ws.Range("D21").Value2 = 0.08429460640631703
$ws.Range("E21").Value2 = 0.06192608065533989
$ws.Range("F21").Value2 = 2.44653680180329
$ws.Range("G21").Value2 = 0.002466953874588532
$ws.Range("K21").Value2 = 3.162693592883954
$ws.Range("M21").Value2 = 0.7679442930468952
$ws.Range("N21").Value2 = 1.771349201146975
$ws.Range("C22").Value2 = 0.08012424006054175
$ws.Range("D22").Value2 = 0.08442762591530339
$ws.Range("E22").Value2 = 0.06530124751426314
$ws.Range("F22").Value2 = 2.512746878677177
$ws.Range("G22").Value2 = 0.002462293948413404
$ws.Range("K22").Value2 = 3.346110539473216
$ws.Range("M22").Value2 = 0.8119702952279511
$ws.Range("N22").Value2 = 1.765589989775876
$ws.Range("C23").Value2 = 0.07739807478677108
$ws.Range("D23").Value2 = 0.08435271265386035
$ws.Range("E23").Value2 = 0.06349655133237775
$ws.Range("F23").Value2 = 2.47727773617305
$ws.Range("G23").Value2 = 0.002464765512415774
$ws.Range("K23").Value2 = 3.24808611687763
$ws.Range("M23").Value2 = 0.7884390162126351
$ws.Range("N23").Value2 = 1.768608542152435
$ws.Range("C24").Value2 = 0.06712685892627235
$ws.Range("D24").Value2 = 0.08415436447454283
$ws.Range("E24").Value2 = 0.05673289345313748
$ws.Range("F24").Value2 = 2.345808255686791
$ws.Range("G24").Value2 = 0.002474468343694697
$ws.Range("K24").Value2 = 2.879678916631519
$ws.Range("M24").Value2 = 0.7000488558952327
$ws.Range("N24").Value2 = 1.781247506934264
$ws.Range("C25").Value2 = 0.05615458789809225
$ws.Range("D25").Value2 = 0.08411838739754529
$ws.Range("E25").Value2 = 0.04958152408476835
$ws.Range("F25").Value2 = 2.209951052228888
$ws.Range("G25").Value2 = 0.002485672290086935
$ws.Range("K25").Value2 = 2.48811465349263
$ws.Range("M25").Value2 = 0.6061974263558199
$ws.Range("N25").Value2 = 1.797414401886257
